{"js": "// Update the date line and the 25 two-digit\u00f7one-digit division problems\n// in the practice-sheet table, per the day's new worksheet content.\n\n// 1) Date heading: \"2025-07-24 Thursday\" -> \"2025-07-25 Friday\"\nconst dateParagraph = context.document.body.paragraphs.getFirst();\ndateParagraph.insertText(\"2025-07-25 Friday\", \"Replace\");\n\n// 2) Table of division problems: replace each cell's text by (row, col)\n// position so identical/duplicate expressions in different cells are each\n// updated independently (text search alone would be ambiguous here).\nconst table = context.document.body.tables.getFirst();\n\n// Each block of 4 table rows holds one row of 5 problems followed by 3\n// blank spacer rows; the populated rows are 0, 4, 8, 12, 16.\nconst newValues = [\n  { row: 0, values: [\"59\u00f74=\", \"84\u00f79=\", \"91\u00f79=\", \"63\u00f77=\", \"40\u00f77=\"] },\n  { row: 4, values: [\"75\u00f76=\", \"39\u00f72=\", \"32\u00f73=\", \"59\u00f74=\", \"29\u00f78=\"] },\n  { row: 8, values: [\"77\u00f76=\", \"85\u00f73=\", \"58\u00f72=\", \"20\u00f79=\", \"28\u00f75=\"] },\n  { row: 12, values: [\"42\u00f74=\", \"29\u00f78=\", \"57\u00f78=\", \"90\u00f74=\", \"22\u00f73=\"] },\n  { row: 16, values: [\"26\u00f75=\", \"57\u00f78=\", \"97\u00f73=\", \"43\u00f77=\", \"52\u00f74=\"] },\n];\n\nfor (const { row, values } of newValues) {\n  for (let col = 0; col < values.length; col++) {\n    table.getCell(row, col).value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 two-digit\u00f7one-digit division problems\n# in the practice-sheet table, per the day's new worksheet content.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading: \"2025-07-24 Thursday\" -> \"2025-07-25 Friday\"\n$d.Paragraphs.Item(1).Range.Text = \"2025-07-25 Friday\"\n\n# 2) Table of division problems: replace each cell's text by (row, col)\n# position so identical/duplicate expressions in different cells are each\n# updated independently (a text-only Find/Replace would be ambiguous here).\n$tbl = $d.Tables.Item(1)\n\n# Each block of 4 table rows holds one row of 5 problems followed by 3\n# blank spacer rows; the populated (1-based) rows are 1, 5, 9, 13, 17.\n$newValues = @{\n    1  = @(\"59\u00f74=\", \"84\u00f79=\", \"91\u00f79=\", \"63\u00f77=\", \"40\u00f77=\")\n    5  = @(\"75\u00f76=\", \"39\u00f72=\", \"32\u00f73=\", \"59\u00f74=\", \"29\u00f78=\")\n    9  = @(\"77\u00f76=\", \"85\u00f73=\", \"58\u00f72=\", \"20\u00f79=\", \"28\u00f75=\")\n    13 = @(\"42\u00f74=\", \"29\u00f78=\", \"57\u00f78=\", \"90\u00f74=\", \"22\u00f73=\")\n    17 = @(\"26\u00f75=\", \"57\u00f78=\", \"97\u00f73=\", \"43\u00f77=\", \"52\u00f74=\")\n}\n\nforeach ($row in $newValues.Keys) {\n    $values = $newValues[$row]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $tbl.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
